$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")
$ws.Visible = -1

$ws.Range("AF1:AF144").Value2 = $ws.Range("AE1:AE144").Value2
$ws.Range("AE1:AE144").Value2 = $ws.Range("AD1:AD144").Value2
$ws.Range("AD1:AD144").Value2 = $ws.Range("AC1:AC144").Value2
$ws.Range("AC1:AC144").Value2 = $ws.Range("AB1:AB144").Value2
$ws.Range("AB1:AB144").Value2 = $ws.Range("AA1:AA144").Value2
$ws.Range("AA1:AA144").Value2 = $ws.Range("Z1:Z144").Value2

Write-Output "AF1=" $ws.Range("AF1").Value2
Write-Output "AF2=" $ws.Range("AF2").Value2
Write-Output "AF27=" $ws.Range("AF27").Value2
Write-Output "AE1=" $ws.Range("AE1").Value2
Write-Output "AE2=" $ws.Range("AE2").Value2
Write-Output "AE8=" $ws.Range("AE8").Value2
Write-Output "AD1=" $ws.Range("AD1").Value2
Write-Output "AD17=" $ws.Range("AD17").Value2
Write-Output "AC1=" $ws.Range("AC1").Value2
Write-Output "AC10=" $ws.Range("AC10").Value2
Write-Output "AB1=" $ws.Range("AB1").Value2
Write-Output "AB8=" $ws.Range("AB8").Value2
Write-Output "AA1=" $ws.Range("AA1").Value2
Write-Output "AA144=" $ws.Range("AA144").Value2
Write-Output "Z1=" $ws.Range("Z1").Value2
Write-Output "Z2=" $ws.Range("Z2").Value2

$ws.Range("Z1:Z144").ClearContents()
$ws.Range("Z1").Value2 = "tn.5250"
$ws.Range("Z2").Value2 = "close(profile)"
$ws.Range("Z3").Value2 = "open(profile)"
$ws.Range("Z4").Value2 = "saveText(profile,var)"
$ws.Range("Z5").Value2 = "typeKeys(profile,keystrokes)"
$ws.Range("Z6").Value2 = "updateScreenFields(profile)"

Write-Output "new Z1=" $ws.Range("Z1").Value2
Write-Output "new Z6=" $ws.Range("Z6").Value2
Write-Output "new Z7=" $ws.Range("Z7").Value2

# Column A (target list): insert tn.5250 between step (A25) and web (A26)
$ws.Range("A32").Value2 = $ws.Range("A31").Value2
$ws.Range("A31").Value2 = $ws.Range("A30").Value2
$ws.Range("A30").Value2 = $ws.Range("A29").Value2
$ws.Range("A29").Value2 = $ws.Range("A28").Value2
$ws.Range("A28").Value2 = $ws.Range("A27").Value2
$ws.Range("A27").Value2 = $ws.Range("A26").Value2
$ws.Range("A26").Value2 = "tn.5250"

Write-Output "A25=" $ws.Range("A25").Value2
Write-Output "A26=" $ws.Range("A26").Value2
Write-Output "A27=" $ws.Range("A27").Value2
Write-Output "A28=" $ws.Range("A28").Value2
Write-Output "A29=" $ws.Range("A29").Value2
Write-Output "A30=" $ws.Range("A30").Value2
Write-Output "A31=" $ws.Range("A31").Value2
Write-Output "A32=" $ws.Range("A32").Value2

# Column K (image commands): insert ocr(image,saveVar) between crop (K5) and resize (K6)
$ws.Range("K8").Value2 = $ws.Range("K7").Value2
$ws.Range("K7").Value2 = $ws.Range("K6").Value2
$ws.Range("K6").Value2 = "ocr(image,saveVar)"
$ws.Range("K2").Value2 = "colorbit(image,bit,saveTo)"

Write-Output "K2=" $ws.Range("K2").Value2
Write-Output "K5=" $ws.Range("K5").Value2
Write-Output "K6=" $ws.Range("K6").Value2
Write-Output "K7=" $ws.Range("K7").Value2
Write-Output "K8=" $ws.Range("K8").Value2

# Update defined names
$wb.Names.Item("image").RefersTo = "='#system'!`$K`$2:`$K`$8"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$144"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"
$wb.Names.Add("tn.5250", "='#system'!`$Z`$2:`$Z`$6")

